# The presentation's colour theme ("Integral", stored as ppt/theme/theme1.xml,
# used by the slide master / all slides / layouts) is replaced with the
# default "Office Theme" colour palette (previously stored only in
# ppt/theme/theme2.xml, which the notes master used).
#
# The font scheme and format scheme are already identical between the two
# themes, so realigning the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) reproduces the effective colour change.

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index -> (theme slot, new "Office Theme" RGB as decimal R+G*256+B*65536)
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
